# Generate Report for Handback
# This script updates the localization-status workbook to reflect that the
# two handed-back files are now "in sync with en-US", and fills in the
# "Latest Target File" / "Latest Handback File" / "Latest Handback DateTime"
# columns (I/J/K) on the per-language sheets, linking the target file name
# the same way the source file name (column A) is already linked.

$wb = $excel.ActiveWorkbook

$statusText = "Handed back: in sync with en-US"

$mdName1 = "ae0f2879-a72e-4765-8ae4-4d529a296e95.md"
$mdName2 = "d773a386-f4b4-44cd-b0fa-74f1a6201a34.md"

$mdUrl1 = "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/557bc2f79f8d0b25847f27859488fac5e110923e/e2e/ae0f2879-a72e-4765-8ae4-4d529a296e95.md"
$mdUrl2 = "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/557bc2f79f8d0b25847f27859488fac5e110923e/e2e/d773a386-f4b4-44cd-b0fa-74f1a6201a34.md"

# Matches the color already used for the existing "Source File Name" hyperlinks
# (font color FF6495ED == RGB(100,149,237), i.e. cornflower blue).
$hyperlinkColor = 15570276

# ----------------------------------------------------------------------
# Overview sheet: flip the per-language status cells to the handback text
# ----------------------------------------------------------------------
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("E2").Value = $statusText
$wsOverview.Range("F2").Value = $statusText
$wsOverview.Range("E3").Value = $statusText
$wsOverview.Range("F3").Value = $statusText

$wsOverview.Columns.Item(5).ColumnWidth = 29.15
$wsOverview.Columns.Item(6).ColumnWidth = 29.15

# ----------------------------------------------------------------------
# Helper to fill in the Latest Target File / Latest Handback File /
# Latest Handback DateTime columns for a language sheet, and link the
# new "Latest Target File" cell the same way column A is linked.
# NOTE: this runtime's PowerShell only reliably supports *positional*
# function arguments, so avoid named (-Param value) calling syntax.
# ----------------------------------------------------------------------
function Update-LanguageSheet($SheetName, $XlfName1, $XlfName2, $HandbackDateTime) {

    $ws = $wb.Worksheets.Item($SheetName)

    # Status column (Status == shared text also used on the Overview sheet)
    $ws.Range("C2").Value = $statusText
    $ws.Range("C3").Value = $statusText

    # Row 2 (ae0f2879... file)
    $i2 = $ws.Range("I2")
    $i2.Value = $mdName1
    $ws.Hyperlinks.Add($i2, $mdUrl1, "", "", $mdName1) | Out-Null
    $i2.Font.Color = $hyperlinkColor
    $i2.Font.Underline = $True

    $ws.Range("J2").Value = $XlfName1
    $ws.Range("K2").Value = $HandbackDateTime

    # Row 3 (d773a386... file)
    $i3 = $ws.Range("I3")
    $i3.Value = $mdName2
    $ws.Hyperlinks.Add($i3, $mdUrl2, "", "", $mdName2) | Out-Null
    $i3.Font.Color = $hyperlinkColor
    $i3.Font.Underline = $True

    $ws.Range("J3").Value = $XlfName2
    $ws.Range("K3").Value = $HandbackDateTime

    # Column widths: Status (C), Latest Target File (I), Latest Handback File (J)
    $ws.Columns.Item(3).ColumnWidth = 29.15
    $ws.Columns.Item(9).ColumnWidth = 39.17
    $ws.Columns.Item(10).ColumnWidth = 39.17
}

Update-LanguageSheet "zh-cn" `
    "ae0f2879-a72e-4765-8ae4-4d529a296e95.028580630c63055c197065015fdac7688747d68f.zh-cn.xlf" `
    "d773a386-f4b4-44cd-b0fa-74f1a6201a34.e63042db4ab5f5808e0db818a521ec6920507d43.zh-cn.xlf" `
    "2016-08-22 12:23:16"

Update-LanguageSheet "de-de" `
    "ae0f2879-a72e-4765-8ae4-4d529a296e95.028580630c63055c197065015fdac7688747d68f.de-de.xlf" `
    "d773a386-f4b4-44cd-b0fa-74f1a6201a34.e63042db4ab5f5808e0db818a521ec6920507d43.de-de.xlf" `
    "2016-08-22 12:23:23"

Write-Host "Handback report generated."
